$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) "Course" -> "Course(s) and number(s) of sections"  (cell D3) ---
$ws.Range("D3").Value = "Course(s) and number(s) of sections"

# --- 2) Drop "UCWU " before "Math Organizing Committee" (cell F7) ---
$ws.Range("F7").Value = "Would you be willing to share course materials with the Math Organizing Committee to support future instructors/TAs/FAs?"

# --- 3) Instructions cell (F10, merged F10:G21) ---
$instrCell = $ws.Range("F10")

$instrText = @"
(1) Please keep track of any time spent on teaching-related work and enter the number of hours spent during each week in the corresponding box in column B.  
—If you only have an estimate, please add an asterisk (*) after the number
—If you are teaching more than one course, you may use one sheet for each course or combine your hours in a single sheet
—To help keep track of time, we suggest using an app like Toggl Track

(2) Please enter other course information in column E and your personal information in column G.  If you are recording multiple courses on this sheet, please indicate for each course as appropriate.

(3) At the end of the term, the Math Organizing Committee will reach out with instructions to submit this spreadsheet.

Please contact Michael Kopreski (michaelkopreski@gmail.com) or any other Math OC members with questions or comments.  Thank you for your time and support!
"@

$instrCell.Value = $instrText

# Re-apply character-level formatting (rich text runs) on top of the new text.
# Offsets (1-based) computed from $instrText above.

# "any" -> italic
$instrCell.Characters(26, 3).Font.Italic = $true

# "If you only have an estimate, please add an asterisk (*) after the number\n" -> bold
$instrCell.Characters(160, 74).Font.Bold = $true

# newline right after the new "one sheet for each course..." sentence -> bold
$instrCell.Characters(354, 1).Font.Bold = $true

# "Toggl Track" -> underline
$instrCell.Characters(413, 11).Font.Underline = $true

# newline right after "Toggl Track" -> bold
$instrCell.Characters(424, 1).Font.Bold = $true

# "michaelkopreski@gmail.com" -> underline
$instrCell.Characters(779, 25).Font.Underline = $true
